$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "placeholder"
$s1.Shapes.Item(1).TextFrame.TextRange.Text = "First slide"

$s3 = $p.Slides.Item(3)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "placeholder"
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Third slide"
